$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 3819
$ws.Range("J62").Value = 3573.75
$ws.Range("L62").Value = 3573.75
$ws.Range("N62").Value = -4821.75
$ws.Range("H65").Value = 3819
$ws.Range("J65").Value = 3573.75
$ws.Range("L65").Value = 17868.75
$ws.Range("N65").Value = -24108.75
$ws.Range("H74").Value = 3718.7827
$ws.Range("I74").Value = 3671.6
$ws.Range("J74").Value = 4033.3333
$ws.Range("K74").Value = 3671.6
$ws.Range("L74").Value = 4033.3333
$ws.Range("M74").Value = -2735.6
$ws.Range("N74").Value = -5905.3333
$ws.Range("H77").Value = 3718.7827
$ws.Range("I77").Value = 3671.6
$ws.Range("J77").Value = 4033.3333
$ws.Range("K77").Value = 18358
$ws.Range("L77").Value = 20166.6665
$ws.Range("M77").Value = -13678
$ws.Range("N77").Value = -29526.6665
$ws.Range("H112").Value = 31251150
$ws.Range("I112").Value = 462.5
$ws.Range("J112").Value = 35715532
$ws.Range("K112").Value = 1387.5
$ws.Range("L112").Value = 107146596
$ws.Range("M112").Value = -279.5
$ws.Range("N112").Value = -107148812
$ws.Range("H137").Value = 22228030
$ws.Range("H138").Value = 2616.849
$ws.Range("I138").Value = 1512.129
$ws.Range("J138").Value = 4173.5
$ws.Range("K138").Value = 4536.387
$ws.Range("L138").Value = 12520.5
$ws.Range("M138").Value = 603.6130000000003
$ws.Range("N138").Value = -22800.5
$ws.Range("H141").Value = 5147.5
$ws.Range("I141").Value = 5147.5
$ws.Range("K141").Value = 15442.5
$ws.Range("M141").Value = -10262.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5583.0166
$ws.Range("I32").Value = 6904.93
$ws.Range("K32").Value = 6904.93
$ws.Range("M32").Value = -6617.93
$ws.Range("H63").Value = 4900
$ws.Range("I63").Value = 2800
$ws.Range("J63").Value = 7000
$ws.Range("K63").Value = 2800
$ws.Range("L63").Value = 7000
$ws.Range("M63").Value = -2114
$ws.Range("N63").Value = -8372
$ws.Range("H66").Value = 4900
$ws.Range("I66").Value = 2800
$ws.Range("J66").Value = 7000
$ws.Range("K66").Value = 14000
$ws.Range("L66").Value = 35000
$ws.Range("M66").Value = -10568
$ws.Range("N66").Value = -41864
$ws.Range("H74").Value = 6300.0527
$ws.Range("J74").Value = 10312.546
$ws.Range("L74").Value = 10312.546
$ws.Range("N74").Value = -12060.546
$ws.Range("H77").Value = 6300.0527
$ws.Range("J77").Value = 10312.546
$ws.Range("L77").Value = 51562.73
$ws.Range("N77").Value = -60298.73
$ws.Range("H102").Value = 1827.5625
$ws.Range("I102").Value = 1729
$ws.Range("K102").Value = 1729
$ws.Range("M102").Value = -107
$ws.Range("H125").Value = 45893.75
$ws.Range("J125").Value = 45893.75
$ws.Range("L125").Value = 45893.75
$ws.Range("N125").Value = -55733.75
$ws.Range("H132").Value = 3189.35
$ws.Range("I132").Value = 3200.6428
$ws.Range("K132").Value = 9601.928400000001
$ws.Range("M132").Value = -7071.928400000001
$ws.Range("H135").Value = 52971.6
$ws.Range("J135").Value = 52971.6
$ws.Range("L135").Value = 52971.6
$ws.Range("N135").Value = -63111.6
$ws.Range("H139").Value = 48514.332
$ws.Range("J139").Value = 48514.332
$ws.Range("L139").Value = 48514.332
$ws.Range("N139").Value = -58794.332

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H81").Value = 27265
$ws.Range("J81").Value = 27265
$ws.Range("L81").Value = 27265
$ws.Range("N81").Value = -29387
$ws.Range("H84").Value = 27265
$ws.Range("J84").Value = 27265
$ws.Range("L84").Value = 81795
$ws.Range("N84").Value = -92403
$ws.Range("H86").Value = 1772.2894
$ws.Range("I86").Value = 1625.6552
$ws.Range("J86").Value = 2244.7778
$ws.Range("K86").Value = 1625.6552
$ws.Range("L86").Value = 2244.7778
$ws.Range("M86").Value = -502.6551999999999
$ws.Range("N86").Value = -4490.7778
$ws.Range("H89").Value = 1772.2894
$ws.Range("I89").Value = 1625.6552
$ws.Range("J89").Value = 2244.7778
$ws.Range("K89").Value = 8128.276
$ws.Range("L89").Value = 11223.889
$ws.Range("M89").Value = -2512.276
$ws.Range("N89").Value = -22455.889
$ws.Range("H99").Value = 919.93335
$ws.Range("I99").Value = 978.8
$ws.Range("J99").Value = 802.2
$ws.Range("K99").Value = 978.8
$ws.Range("L99").Value = 802.2
$ws.Range("M99").Value = 519.2
$ws.Range("N99").Value = -3798.2
$ws.Range("H120").Value = 39999
$ws.Range("J120").Value = 39999
$ws.Range("L120").Value = 39999
$ws.Range("H125").Value = 52320
$ws.Range("J125").Value = 52320
$ws.Range("L125").Value = 52320
$ws.Range("N125").Value = -62160
$ws.Range("H134").Value = 75360.266
$ws.Range("I134").Value = 75360.266
$ws.Range("K134").Value = 226080.798
$ws.Range("M134").Value = -223545.798
$ws.Range("H135").Value = 39593.332
$ws.Range("J135").Value = 39593.332
$ws.Range("L135").Value = 39593.332
$ws.Range("N135").Value = -49733.332
$ws.Range("N120").Value = -49675

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3472.5925
$ws.Range("I31").Value = 1411.4286
$ws.Range("J31").Value = 5692.3076
$ws.Range("K31").Value = 1411.4286
$ws.Range("L31").Value = 5692.3076
$ws.Range("M31").Value = -1116.4286
$ws.Range("N31").Value = -6282.3076
$ws.Range("H34").Value = 3472.5925
$ws.Range("I34").Value = 1411.4286
$ws.Range("J34").Value = 5692.3076
$ws.Range("K34").Value = 1411.4286
$ws.Range("L34").Value = 5692.3076
$ws.Range("M34").Value = -1209.4286
$ws.Range("N34").Value = -6096.3076

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H10").Value = 2225
$ws.Range("I10").Value = 1966.6666
$ws.Range("K10").Value = 1966.6666
$ws.Range("M10").Value = -1797.6666
$ws.Range("H132").Value = 3204.9412
$ws.Range("I132").Value = 2336.9167
$ws.Range("K132").Value = 7010.750100000001
$ws.Range("M132").Value = -4480.750100000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1530.8462
$ws.Range("I68").Value = 1350
$ws.Range("J68").Value = 2133.6667
$ws.Range("K68").Value = 1350
$ws.Range("L68").Value = 2133.6667
$ws.Range("M68").Value = -601
$ws.Range("N68").Value = -3631.6667
$ws.Range("H71").Value = 1530.8462
$ws.Range("I71").Value = 1350
$ws.Range("J71").Value = 2133.6667
$ws.Range("K71").Value = 6750
$ws.Range("L71").Value = 10668.3335
$ws.Range("M71").Value = -3006
$ws.Range("N71").Value = -18156.3335
$ws.Range("H122").Value = 22317.6
$ws.Range("I122").Value = 27072
$ws.Range("K122").Value = 81216
$ws.Range("M122").Value = -78766
$ws.Range("H132").Value = 11793.8
$ws.Range("I132").Value = 16434.223
$ws.Range("J132").Value = 4833.1665
$ws.Range("K132").Value = 49302.66900000001
$ws.Range("L132").Value = 14499.4995
$ws.Range("M132").Value = -46772.66900000001
$ws.Range("N132").Value = -19559.4995
$ws.Range("H136").Value = 2306.6667
$ws.Range("J136").Value = 2828.5715
$ws.Range("L136").Value = 8485.7145
$ws.Range("N136").Value = -13585.7145

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 10492.2
$ws.Range("I62").Value = 8584.4
$ws.Range("J62").Value = 12400
$ws.Range("K62").Value = 8584.4
$ws.Range("L62").Value = 12400
$ws.Range("M62").Value = -7960.4
$ws.Range("N62").Value = -13648
$ws.Range("H65").Value = 10492.2
$ws.Range("I65").Value = 8584.4
$ws.Range("J65").Value = 12400
$ws.Range("K65").Value = 42922
$ws.Range("L65").Value = 62000
$ws.Range("M65").Value = -39802
$ws.Range("N65").Value = -68240
